# "Generate Report for Handback" re-run:
# the eaebf39c-c09c-44b7-9f1e-18c264aef4fe.* entry has dropped out of the
# report (handback finished / entry retired) and the remaining
# 39b07019-896a-4d16-842b-bb42829f0703.* entry's handoff/handback
# timestamps were refreshed to a later run.
#
# Concretely, on every sheet the 2nd data row (row 3) - which held the
# eaebf39c... file - is removed (with its hyperlinks), and on the
# per-language sheets the "Correspond Handoff/Handback DateTime" cells
# for the remaining row are bumped to the new timestamps.

$wb = $excel.ActiveWorkbook

function Remove-HyperlinkAt($ws, $addr) {
    # NOTE: $range.Hyperlinks.Delete() removes *every* hyperlink on the
    # worksheet in this host, not just the ones touching $range, so find
    # the specific Hyperlink object whose Range matches and delete that.
    $toDelete = $null
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $toDelete = $hl
        }
    }
    if ($toDelete -ne $null) {
        $toDelete.Delete()
    }
}

# ---- Overview sheet: drop row 3 (eaebf39c... entry) ----
$ws1 = $wb.Worksheets.Item("Overview")
Remove-HyperlinkAt $ws1 "`$A`$3"
$ws1.Rows.Item(3).Delete()

# ---- zh-cn sheet: refresh row 2 timestamps, then drop row 3 ----
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("E2").Value = "2016-03-13 06:51:29"
$ws2.Range("H2").Value = "2016-03-13 06:51:51"
Remove-HyperlinkAt $ws2 "`$A`$3"
Remove-HyperlinkAt $ws2 "`$B`$3"
Remove-HyperlinkAt $ws2 "`$D`$3"
Remove-HyperlinkAt $ws2 "`$F`$3"
Remove-HyperlinkAt $ws2 "`$G`$3"
$ws2.Rows.Item(3).Delete()

# ---- de-de sheet: refresh row 2 timestamps, then drop row 3 ----
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("E2").Value = "2016-03-13 06:51:32"
$ws3.Range("H2").Value = "2016-03-13 06:51:57"
Remove-HyperlinkAt $ws3 "`$A`$3"
Remove-HyperlinkAt $ws3 "`$B`$3"
Remove-HyperlinkAt $ws3 "`$D`$3"
Remove-HyperlinkAt $ws3 "`$F`$3"
Remove-HyperlinkAt $ws3 "`$G`$3"
$ws3.Rows.Item(3).Delete()
